$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 895.25
$ws.Range("I18").Value = 893.6667
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 893.6667
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = -609.6667
$ws.Range("N18").Value = -1468
$ws.Range("H19").Value = 638.34784
$ws.Range("I19").Value = 722.1
$ws.Range("J19").Value = 573.9231
$ws.Range("K19").Value = 722.1
$ws.Range("L19").Value = 573.9231
$ws.Range("M19").Value = -547.1
$ws.Range("N19").Value = -923.9231
$ws.Range("H70").Value = 12445740
$ws.Range("I70").Value = 24890324
$ws.Range("J70").Value = 1156.6666
$ws.Range("K70").Value = 74670972
$ws.Range("L70").Value = 3469.9998
$ws.Range("M70").Value = -74670702
$ws.Range("N70").Value = -4009.9998
$ws.Range("H73").Value = 12445740
$ws.Range("I73").Value = 24890324
$ws.Range("J73").Value = 1156.6666
$ws.Range("K73").Value = 74670972
$ws.Range("L73").Value = 3469.9998
$ws.Range("M73").Value = -74670036
$ws.Range("N73").Value = -5341.9998
$ws.Range("H112").Value = 1940.8
$ws.Range("I112").Value = 525
$ws.Range("J112").Value = 2063.913
$ws.Range("K112").Value = 1575
$ws.Range("L112").Value = 6191.739
$ws.Range("M112").Value = -467
$ws.Range("N112").Value = -8407.739
$ws.Range("H132").Value = 29457516
$ws.Range("I132").Value = 30071194
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 90213582
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -90211052
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 2393983.5
$ws.Range("I137").Value = 855.1739
$ws.Range("J137").Value = 5290928.5
$ws.Range("K137").Value = 2565.5217
$ws.Range("L137").Value = 15872785.5
$ws.Range("M137").Value = -15.52170000000024
$ws.Range("N137").Value = -15877885.5
$ws.Range("H138").Value = 1697.6323
$ws.Range("I138").Value = 657.7273
$ws.Range("J138").Value = 2194.9783
$ws.Range("K138").Value = 1973.1819
$ws.Range("L138").Value = 6584.9349
$ws.Range("M138").Value = 3166.8181
$ws.Range("N138").Value = -16864.9349
$ws.Range("H141").Value = 2866.8857
$ws.Range("I141").Value = 2001.0333
$ws.Range("J141").Value = 8062
$ws.Range("K141").Value = 6003.0999
$ws.Range("L141").Value = 24186
$ws.Range("M141").Value = -823.0999000000002
$ws.Range("N141").Value = -34546
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 476368.62
$ws.Range("I32").Value = 3486.3462
$ws.Range("J32").Value = 2525525.2
$ws.Range("K32").Value = 3486.3462
$ws.Range("L32").Value = 2525525.2
$ws.Range("M32").Value = -3199.3462
$ws.Range("N32").Value = -2526099.2
$ws.Range("H132").Value = 3908524.2
$ws.Range("I132").Value = 8334485.5
$ws.Range("J132").Value = 3264.4707
$ws.Range("K132").Value = 25003456.5
$ws.Range("L132").Value = 9793.4121
$ws.Range("M132").Value = -25000926.5
$ws.Range("N132").Value = -14853.4121
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7588872
$ws.Range("I134").Value = 8561489
$ws.Range("J134").Value = 2462.8
$ws.Range("K134").Value = 25684467
$ws.Range("L134").Value = 7388.400000000001
$ws.Range("M134").Value = -25681932
$ws.Range("N134").Value = -12458.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11910223
$ws.Range("I132").Value = 15152965
$ws.Range("J132").Value = 20168.834
$ws.Range("K132").Value = 45458895
$ws.Range("L132").Value = 60506.50199999999
$ws.Range("M132").Value = -45456365
$ws.Range("N132").Value = -65566.50199999999
$ws.Range("H134").Value = 22630782
$ws.Range("I134").Value = 29763080
$ws.Range("J134").Value = 3908500
$ws.Range("K134").Value = 89289240
$ws.Range("L134").Value = 11725500
$ws.Range("M134").Value = -89286705
$ws.Range("N134").Value = -11730570
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1004.5833
$ws.Range("I2").Value = 278
$ws.Range("J2").Value = 1367.875
$ws.Range("K2").Value = 1668
$ws.Range("L2").Value = 8207.25
$ws.Range("M2").Value = -1555
$ws.Range("N2").Value = -8433.25
$ws.Range("H11").Value = 2347.7778
$ws.Range("I11").Value = 1810
$ws.Range("J11").Value = 3423.3333
$ws.Range("K11").Value = 5430
$ws.Range("L11").Value = 10269.9999
$ws.Range("M11").Value = -5290
$ws.Range("N11").Value = -10549.9999
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 300
$ws.Range("K16").Value = 900
$ws.Range("M16").Value = -727
$ws.Range("H20").Value = 1900
$ws.Range("J20").Value = 1880
$ws.Range("L20").Value = 5640
$ws.Range("N20").Value = -6094
$ws.Range("H26").Value = 719.53845
$ws.Range("I26").Value = 90.8
$ws.Range("J26").Value = 1112.5
$ws.Range("K26").Value = 272.4
$ws.Range("L26").Value = 3337.5
$ws.Range("M26").Value = 15.60000000000002
$ws.Range("N26").Value = -3913.5
$ws.Range("H63").Value = 5976.5557
$ws.Range("I63").Value = 2226
$ws.Range("J63").Value = 7048.143
$ws.Range("K63").Value = 6678
$ws.Range("L63").Value = 21144.429
$ws.Range("M63").Value = -5929
$ws.Range("N63").Value = -22642.429
$ws.Range("H66").Value = 5976.5557
$ws.Range("I66").Value = 2226
$ws.Range("J66").Value = 7048.143
$ws.Range("K66").Value = 20034
$ws.Range("L66").Value = 63433.287
$ws.Range("M66").Value = -16290
$ws.Range("N66").Value = -70921.287
$ws.Range("H68").Value = 228281.3
$ws.Range("I68").Value = 509.55554
$ws.Range("J68").Value = 529743.9
$ws.Range("K68").Value = 1528.66662
$ws.Range("L68").Value = 1589231.7
$ws.Range("M68").Value = -717.66662
$ws.Range("N68").Value = -1590853.7
$ws.Range("H71").Value = 228281.3
$ws.Range("I71").Value = 509.55554
$ws.Range("J71").Value = 529743.9
$ws.Range("K71").Value = 4585.99986
$ws.Range("L71").Value = 4767695.100000001
$ws.Range("M71").Value = -529.9998599999999
$ws.Range("N71").Value = -4775807.100000001
$ws.Range("H75").Value = 83333704
$ws.Range("I75").Value = 492
$ws.Range("J75").Value = 333333340
$ws.Range("K75").Value = 1476
$ws.Range("L75").Value = 1000000020
$ws.Range("M75").Value = -478
$ws.Range("N75").Value = -1000002016
$ws.Range("H78").Value = 83333704
$ws.Range("I78").Value = 492
$ws.Range("J78").Value = 333333340
$ws.Range("K78").Value = 4428
$ws.Range("L78").Value = 3000000060
$ws.Range("M78").Value = 564
$ws.Range("N78").Value = -3000010044
$ws.Range("H107").Value = 650.4337
$ws.Range("I107").Value = 322.64178
$ws.Range("J107").Value = 2023.0625
$ws.Range("K107").Value = 967.92534
$ws.Range("L107").Value = 6069.1875
$ws.Range("M107").Value = 952.07466
$ws.Range("N107").Value = -9909.1875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 16625
$ws.Range("J95").Value = 16625
$ws.Range("L95").Value = 16625
$ws.Range("N95").Value = -22117
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 25003178
$ws.Range("I122").Value = 693.3333
$ws.Range("J122").Value = 35718530
$ws.Range("K122").Value = 2079.9999
$ws.Range("L122").Value = 107155590
$ws.Range("M122").Value = 370.0001000000002
$ws.Range("N122").Value = -107160490
$ws.Range("H132").Value = 2447701.2
$ws.Range("I132").Value = 2365047.5
$ws.Range("J132").Value = 2676214.8
$ws.Range("K132").Value = 7095142.5
$ws.Range("L132").Value = 8028644.399999999
$ws.Range("M132").Value = -7092612.5
$ws.Range("N132").Value = -8033704.399999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22900930
$ws.Range("I132").Value = 27908258
$ws.Range("J132").Value = 1369418.2
$ws.Range("K132").Value = 83724774
$ws.Range("L132").Value = 4108254.6
$ws.Range("M132").Value = -83722244
$ws.Range("N132").Value = -4113314.6
